# Updated cryptos list on Mon Feb 27 05:42:59 UTC 2023 with GitHub Actions
#
# Refreshes the Price (D) / Volume(1h) (E) columns for each coin row, and
# swaps the Aptos / EnergySwap rows (43 <-> 44) back into rank order with
# their refreshed figures.
#
# Price values are plain text (e.g. "1.001", "23.479.39") in the source
# sheet. Excel's Range.Value setter auto-coerces anything that parses as a
# plain number (losing formatting like trailing zeros: "1.000" -> 1,
# "151.40" -> 151.4), so for any new price that looks numeric we
# temporarily force the cell to Text format, assign it, then
# ClearFormats() to drop the temporary number-format style again (keeping
# the cell's style untouched, same as the rest of the sheet). Values like
# "23.516.02" (two dots) never round-trip as a number, so they can be
# assigned directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PriceValue {
    param($range, [string]$value)
    if ($value -match '^-?[0-9]+(\.[0-9]+)?$') {
        $range.NumberFormat = "@"
        $range.Value = $value
        $range.ClearFormats()
    } else {
        $range.Value = $value
    }
}

function Set-VolumeValue {
    param($range, [string]$pct)
    $range.Value = "  $pct  "
}

# row -> new price, new volume%
$updates = [ordered]@{
    2  = @("23.516.02",  "+1.37%")
    3  = @("1.640.15",   "+2.35%")
    4  = @("1.001",      "-0.02%")
    5  = @("307.62",     "+1.53%")
    6  = @("1.001",      "+0.06%")
    7  = @("0.3774",     "-0.26%")
    8  = @("52.28",      "+0.59%")
    9  = @("0.3651",     "+0.99%")
    10 = @("1.271",      "+0.40%")
    11 = @("0.08178",    "+0.55%")
    13 = @("23.06",      "+1.83%")
    14 = @("6.651",      "+1.02%")
    15 = @("0.00001281", "+2.63%")
    16 = @("7.392",      "+0.02%")
    17 = @("1.643.17",   "+2.67%")
    18 = @("94.86",      "+1.07%")
    20 = @("18.25",      "+1.01%")
    21 = @("6.555",      "+0.11%")
    22 = @("1.000",      "-0.05%")
    23 = @("23.508.98",  "+1.37%")
    24 = @("12.83",      "-0.84%")
    25 = @("3.110",      "+4.26%")
    26 = @("2.419",      "+1.27%")
    27 = @("21.29",      "+0.41%")
    28 = @("151.40",     "+1.68%")
    29 = @("5.359",      "+2.09%")
    30 = @("135.55",     "+1.29%")
    31 = @("2.356",      "-0.81%")
    32 = @("1.817.95",   "+2.27%")
    33 = @("6.791",      "-0.66%")
    34 = @("0.9664",     "-0.53%")
    35 = @("0.02835",    "+4.37%")
    36 = @("10.37",      "+0.33%")
    37 = @("0.07369",    "-2.05%")
    38 = @("0.2539",     "+1.14%")
    39 = @("6.188",      "+1.13%")
    40 = @("0.08861",    "+0.58%")
    41 = @("1.386",      "+1.69%")
    42 = @("0.7123",     "+0.38%")
    45 = @("0.6564",     "+0.44%")
    46 = @("2.345",      "+1.41%")
    48 = @("4.033",      "+0.59%")
    49 = @("0.07969",    "+0.15%")
    50 = @("129.48",     "-1.97%")
    51 = @("1.212",      "+0.77%")
}

# rows with only a Volume(1h) refresh (price column untouched in the diff)
$volumeOnlyUpdates = [ordered]@{
    12 = "-0.02%"
    47 = "+0.04%"
}

foreach ($row in $updates.Keys) {
    $pair = $updates[$row]
    Set-PriceValue $ws.Range("D$row") $pair[0]
    Set-VolumeValue $ws.Range("E$row") $pair[1]
}

foreach ($row in $volumeOnlyUpdates.Keys) {
    Set-VolumeValue $ws.Range("E$row") $volumeOnlyUpdates[$row]
}

# Rows 43/44 swap rank position: EnergySwap (was #43) and Aptos (was #44)
# trade places, each carrying its own refreshed price/volume.
$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-PriceValue $ws.Range("D43") "12.53"
Set-VolumeValue $ws.Range("E43") "+0.02%"

$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-PriceValue $ws.Range("D44") "16.17"
Set-VolumeValue $ws.Range("E44") "+3.96%"
